$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G9").Value = 1.0
$ws.Range("G11").Value = 1.0
$ws.Range("G17").Value = 4.0
$ws.Range("H17").Value = 4.0
$ws.Range("J17").Value = "Yes"
$ws.Range("K30").Value = "Yes"
$ws.Range("K38").Value = "Yes"
$ws.Range("G43").Value = 4.0
$ws.Range("H43").Value = 4.0
$ws.Range("J43").Value = "Yes"
$ws.Range("H52").Value = 1.0
$ws.Range("H53").Value = 1.0
$ws.Range("G54").Value = 3.0
$ws.Range("H54").Value = 3.0
$ws.Range("K63").Value = "Yes"
$ws.Range("K80").Value = "Yes"
$ws.Range("I87").Value = 1.0
$ws.Range("K91").Value = "Yes"
$ws.Range("K109").Value = "Yes"
$ws.Range("H111").Value = 4.0
$ws.Range("J111").Value = "Yes"
$ws.Range("K117").Value = "Yes"
$ws.Range("G122").Value = 1.0
$ws.Range("K122").Value = "Yes"
$ws.Range("K126").Value = "Yes"
$ws.Range("G143").Value = 1.0
$ws.Range("K143").Value = "Yes"
$ws.Range("G145").Value = 4.0
$ws.Range("H145").Value = 3.0
$ws.Range("K154").Value = "Yes"
$ws.Range("K160").Value = "Yes"
$ws.Range("K174").Value = "Yes"
